# UserListDataSet.xlsx update:
#  - Remove the PurchaserDetails, Menu, and Categories sheets (no longer used by the tests).
#  - Update the SignUp sheet's sample username/password row from "chakk27" to "chakk35".

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets("SignUp")
$ws.Range("A2").Value = "chakk35"
$ws.Range("B2").Value = "chakk35"

$wb.Worksheets("PurchaserDetails").Delete()
$wb.Worksheets("Menu").Delete()
$wb.Worksheets("Categories").Delete()

# Leave the SignUp tab as the selected/active sheet, matching the saved workbook state.
$ws.Select()
